$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A2:B23) in descending order by the value column (B)
$range = $ws.Range("A2:B23")
$sortKey = $ws.Range("B2:B23")
$range.Sort($sortKey, 2)

# After sorting descending by value, the zero-value rows (Russian, Uzbek)
# end up at the bottom (rows 22 and 23). Remove them.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
